$d = $word.ActiveDocument

# Locate the old "Fecha: 9 septiembre 2020" run and replace its text with
# the corrected date. This keeps it as a single run for the moment.
$old = "Fecha: 9 septiembre 2020"
$new = "Fecha: 22 febrero 2021"

$find = $d.Content
$found = $find.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target text '$old'"
}

# Re-locate the replaced text so we know its exact Start/End.
$target = $d.Content
$target.Find.Execute($new, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $target.Start
$end = $target.End

# The committed version splits the single run into six runs with identical
# formatting (only <w:lang w:val="es-CO"/>):
#   "Fecha: " | "22" | " " | "febrero" | " 202" | "1"
# Word naturally merges adjacent runs that share identical formatting, so a
# plain InsertAfter/Text assignment collapses back into one run. Toggling a
# character property (Bold on, then back off) on a sub-range forces Word to
# seat that sub-range as its own run without altering the visible
# formatting, which is how we reproduce the exact run layout.
$boundaries = @(7, 9, 10, 17, 21)
foreach ($b in $boundaries) {
    $sub = $d.Range($start + $b, $end)
    $sub.Bold = 1
    $sub.Bold = 0
}
